$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-04-27 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-28 Sunday", 2) | Out-Null
$d.Content.Find.Execute("169÷7=24, 1", $true, $false, $false, $false, $false, $true, 1, $false, "546÷4=136, 2", 2) | Out-Null
$d.Content.Find.Execute("455÷6=75, 5", $true, $false, $false, $false, $false, $true, 1, $false, "384÷9=42, 6", 2) | Out-Null
$d.Content.Find.Execute("628÷3=209, 1", $true, $false, $false, $false, $false, $true, 1, $false, "869÷5=173, 4", 2) | Out-Null
$d.Content.Find.Execute("387÷3=129, 0", $true, $false, $false, $false, $false, $true, 1, $false, "205÷7=29, 2", 2) | Out-Null
$d.Content.Find.Execute("879÷6=146, 3", $true, $false, $false, $false, $false, $true, 1, $false, "371÷8=46, 3", 2) | Out-Null
$d.Content.Find.Execute("679÷6=113, 1", $true, $false, $false, $false, $false, $true, 1, $false, "437÷8=54, 5", 2) | Out-Null
$d.Content.Find.Execute("892÷8=111, 4", $true, $false, $false, $false, $false, $true, 1, $false, "656÷7=93, 5", 2) | Out-Null
$d.Content.Find.Execute("518÷2=259, 0", $true, $false, $false, $false, $false, $true, 1, $false, "130÷4=32, 2", 2) | Out-Null
$d.Content.Find.Execute("635÷6=105, 5", $true, $false, $false, $false, $false, $true, 1, $false, "447÷9=49, 6", 2) | Out-Null
$d.Content.Find.Execute("793÷2=396, 1", $true, $false, $false, $false, $false, $true, 1, $false, "934÷7=133, 3", 2) | Out-Null
$d.Content.Find.Execute("511÷2=255, 1", $true, $false, $false, $false, $false, $true, 1, $false, "554÷3=184, 2", 2) | Out-Null
$d.Content.Find.Execute("379÷5=75, 4", $true, $false, $false, $false, $false, $true, 1, $false, "895÷7=127, 6", 2) | Out-Null
$d.Content.Find.Execute("305÷5=61, 0", $true, $false, $false, $false, $false, $true, 1, $false, "998÷5=199, 3", 2) | Out-Null
$d.Content.Find.Execute("182÷2=91, 0", $true, $false, $false, $false, $false, $true, 1, $false, "473÷4=118, 1", 2) | Out-Null
$d.Content.Find.Execute("213÷4=53, 1", $true, $false, $false, $false, $false, $true, 1, $false, "269÷8=33, 5", 2) | Out-Null
$d.Content.Find.Execute("184÷9=20, 4", $true, $false, $false, $false, $false, $true, 1, $false, "138÷6=23, 0", 2) | Out-Null
$d.Content.Find.Execute("647÷8=80, 7", $true, $false, $false, $false, $false, $true, 1, $false, "675÷5=135, 0", 2) | Out-Null
$d.Content.Find.Execute("483÷4=120, 3", $true, $false, $false, $false, $false, $true, 1, $false, "468÷7=66, 6", 2) | Out-Null
$d.Content.Find.Execute("790÷5=158, 0", $true, $false, $false, $false, $false, $true, 1, $false, "750÷6=125, 0", 2) | Out-Null
$d.Content.Find.Execute("680÷8=85, 0", $true, $false, $false, $false, $false, $true, 1, $false, "400÷7=57, 1", 2) | Out-Null
$d.Content.Find.Execute("976÷7=139, 3", $true, $false, $false, $false, $false, $true, 1, $false, "424÷9=47, 1", 2) | Out-Null
$d.Content.Find.Execute("932÷2=466, 0", $true, $false, $false, $false, $false, $true, 1, $false, "306÷7=43, 5", 2) | Out-Null
$d.Content.Find.Execute("542÷9=60, 2", $true, $false, $false, $false, $false, $true, 1, $false, "366÷2=183, 0", 2) | Out-Null
$d.Content.Find.Execute("504÷6=84, 0", $true, $false, $false, $false, $false, $true, 1, $false, "176÷5=35, 1", 2) | Out-Null
$d.Content.Find.Execute("435÷6=72, 3", $true, $false, $false, $false, $false, $true, 1, $false, "543÷7=77, 4", 2) | Out-Null
